$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows: MEDIA (102), DESVIO PADRAO (103), TAXA DE SUCESSO (104) ---

# Labels (column A)
$ws.Range("A102").Value = "MEDIA"
$ws.Range("A103").Value = "DESVIO PADRAO"
$ws.Range("A104").Value = "TAXA DE SUCESSO"

# Formulas (columns B:D)
$ws.Range("B102").Formula = "=MEDIAN(B2:B101)"
$ws.Range("C102").Formula = "=MEDIAN(C2:C101)"
$ws.Range("D102").Formula = "=MEDIAN(D2:D101)"

$ws.Range("B103").Formula = "=STDEV.P(B2:B101)"
$ws.Range("C103").Formula = "=STDEV.P(C2:C101)"
$ws.Range("D103").Formula = "=STDEV.P(D2:D101)"

$ws.Range("B104").Formula = '=COUNTIF(B2:B101,">=4")'
$ws.Range("C104").Formula = '=COUNTIF(C2:C101,">=4")'
$ws.Range("D104").Formula = '=COUNTIF(D2:D101,">=4")'

# --- Styling ---

# Label cells: bold white text on dark fill, left/center aligned
foreach ($addr in @("A102", "A103", "A104")) {
  $cell = $ws.Range($addr)
  $cell.Font.Bold = $true
  $cell.Font.ThemeColor = 2
  $cell.Interior.ThemeColor = 1
  $cell.HorizontalAlignment = -4131
  $cell.VerticalAlignment = -4108
}

# Value cells: comma / thousands ("Accounting"-style Comma) number style
$values = $ws.Range("B102:D104")
$values.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'

# Column A widened slightly to fit the new "TAXA DE SUCESSO" label
$ws.Columns.Item(1).ColumnWidth = 17

# --- View / selection ---
[void]$ws.Range("E12").Select()

